$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute("bool, double, byte, char,", $true, $false, $false, $false, $false, $true, 1, $false, "bool, double, char,", 2)
